$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 448, shifting existing rows 448:520 down to 449:521
$ws.Rows.Item(448).Insert()

# Populate the newly inserted row 448 with data
$ws.Cells.Item(448, 1).Value = 9
$ws.Cells.Item(448, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(448, 3).Value = "Metropolitana"
$ws.Cells.Item(448, 4).Value = 44951
$ws.Cells.Item(448, 5).Value = 13
$ws.Cells.Item(448, 6).Value = 100112039
$ws.Cells.Item(448, 7).Value = "Ciboulette"
$ws.Cells.Item(448, 8).Value = "Sin especificar"
$ws.Cells.Item(448, 9).Value = "Primera"
$ws.Cells.Item(448, 10).Value = 340
$ws.Cells.Item(448, 11).Value = 1000
$ws.Cells.Item(448, 12).Value = 1000
$ws.Cells.Item(448, 13).Value = 1000
$ws.Cells.Item(448, 14).Value = "`$/docena de atados"
$ws.Cells.Item(448, 15).Value = "Región Metropolitana"
$ws.Cells.Item(448, 16).Value = 333
$ws.Cells.Item(448, 17).Value = 3
$ws.Cells.Item(448, 18).Value = "Hortaliza"
